# Append a new data row (row 96) to the CityResaleNum sheet, mirroring the
# existing rows: columns A-D hold plain text (date/time/weekday/week are
# stored as literal strings, not native Excel date/time values) while
# columns E-T hold numbers.
#
# A literal apostrophe prefix forces Excel to keep text-looking-like-a-date
# (or text-looking-like-a-number) as text instead of auto-coercing it; the
# style is then copied from the row directly above so the new cells don't
# pick up a stray "quote prefix" number format and stay on the sheet's
# default (unstyled) cell format, matching the rest of the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96
$prevRow = 95

$textValues = @{
    1 = "2023-06-30"
    2 = "22:10:24"
    3 = "Friday"
    4 = "26"
}

foreach ($col in 1..4) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $textValues[$col]
    $cell.Style = $ws.Cells.Item($prevRow, $col).Style
}

$numericValues = @{
    5  = 123466
    6  = 134617
    7  = 163772
    8  = 133779
    9  = 176808
    10 = 115793
    11 = 204998
    12 = 226224
    13 = 176740
    14 = 104717
    15 = 39818
    16 = 33686
    17 = 52663
    18 = -1
    19 = 36336
    20 = -1
}

foreach ($col in 5..20) {
    $ws.Cells.Item($row, $col).Value = $numericValues[$col]
}
